# Append 7 new match rows (rows 9-15) to the "Abhishek Sharma" sheet,
# mirroring the structure/content of the existing rows 2-8.
#
# All cells in this sheet are stored as text (even numeric-looking values
# such as run counts / strike rates), matching the original scraped data
# convention (numberStoredAsText). We force the new cells to stay text by
# applying a Text number format before writing the values, so numeric-looking
# strings like "31" or "129.16" are not coerced into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ r = 9;  A = " Dubai (DSC)"; B = " October 02 2020";    C = "Sunrisers won by 7 runs";                               D = "Sunrisers Hyderabad"; E = "Chennai Super Kings";         F = "Abhishek Sharma "; G = "31"; H = "24"; I = "4"; J = "1"; K = "129.16" },
    @{ r = 10; A = " Abu Dhabi";   B = " September 29 2020";  C = "Sunrisers won by 15 runs";                              D = "Sunrisers Hyderabad"; E = "Delhi Capitals";             F = "Abhishek Sharma "; G = "1";  H = "1";  I = "0"; J = "0"; K = "100.00" },
    @{ r = 11; A = " Sharjah";     B = " October 04 2020";    C = "Mumbai won by 34 runs";                                 D = "Sunrisers Hyderabad"; E = "Mumbai Indians";             F = "Abhishek Sharma "; G = "10"; H = "13"; I = "0"; J = "0"; K = "76.92" },
    @{ r = 12; A = " Abu Dhabi";   B = " September 26 2020";  C = "KKR won by 7 wickets (with 12 balls remaining)";        D = "Sunrisers Hyderabad"; E = "Kolkata Knight Riders";      F = "Abhishek Sharma "; G = "2";  H = "3";  I = "0"; J = "0"; K = "66.66" },
    @{ r = 13; A = " Sharjah";     B = " October 31 2020";    C = "Sunrisers won by 5 wickets (with 35 balls remaining)";  D = "Sunrisers Hyderabad"; E = "Royal Challengers Bangalore"; F = "Abhishek Sharma "; G = "8";  H = "5";  I = "0"; J = "1"; K = "160.00" },
    @{ r = 14; A = " Dubai (DSC)"; B = " September 21 2020";  C = "RCB won by 10 runs";                                    D = "Sunrisers Hyderabad"; E = "Royal Challengers Bangalore"; F = "Abhishek Sharma "; G = "7";  H = "4";  I = "1"; J = "0"; K = "175.00" },
    @{ r = 15; A = " Dubai (DSC)"; B = " October 08 2020";    C = "Sunrisers won by 69 runs";                              D = "Sunrisers Hyderabad"; E = "Kings XI Punjab";            F = "Abhishek Sharma "; G = "12"; H = "6";  I = "1"; J = "1"; K = "200.00" }
)

foreach ($row in $rows) {
    $r = $row.r

    # G:K hold numeric-looking values (runs/balls/4s/6s/strike-rate). Force
    # those to a Text number format *before* assigning, so they stay text
    # cells (matching the rest of the sheet) instead of being coerced into
    # real numbers. A:F are never numeric-looking, so plain assignment
    # already keeps them as text.
    $ws.Range("G$r`:K$r").NumberFormat = "@"

    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
}
